# Update timestamps as part of "Generate Report for Handback"
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-13 13:40:39"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-13 13:40:28"
$wsZhCn.Range("K2").Value = "2016-10-13 13:41:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-13 13:40:39"
$wsDeDe.Range("K2").Value = "2016-10-13 13:41:29"
